$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force Price (D) and Volume(1h) (E) columns to remain plain text so that
# numeric-looking strings (e.g. "4.24", "521.01") are not auto-converted to numbers.
$textRange = $ws.Range("D2:E51")
$textRange.NumberFormat = "@"

$ws.Range("D2").Value = '57.409.14'
$ws.Range("E2").Value = '  +1.11%  '
$ws.Range("D3").Value = '2.361.71'
$ws.Range("E3").Value = '  +1.53%  '
$ws.Range("E4").Value = '  -0.21%  '
$ws.Range("D5").Value = '521.01'
$ws.Range("E5").Value = '  -0.40%  '
$ws.Range("D6").Value = '135.89'
$ws.Range("E6").Value = '  +0.43%  '
$ws.Range("D7").Value = '0.996'
$ws.Range("E7").Value = '  +0.07%  '
$ws.Range("D8").Value = '0.541'
$ws.Range("E8").Value = '  +0.65%  '
$ws.Range("E9").Value = '  -1.27%  '
$ws.Range("E10").Value = '  +4.97%  '
$ws.Range("E11").Value = '  -0.81%  '
$ws.Range("E12").Value = '  -0.98%  '
$ws.Range("D13").Value = '24.39'
$ws.Range("E13").Value = '  +0.91%  '
$ws.Range("D14").Value = '2.784.80'
$ws.Range("E14").Value = '  +1.85%  '
$ws.Range("D15").Value = '57.396.76'
$ws.Range("E15").Value = '  +0.94%  '
$ws.Range("E16").Value = '  -0.25%  '
$ws.Range("D17").Value = '2.372.82'
$ws.Range("E17").Value = '  +2.18%  '
$ws.Range("E18").Value = '  +0.28%  '
$ws.Range("D19").Value = '330.16'
$ws.Range("E19").Value = '  +2.31%  '
$ws.Range("D20").Value = '4.24'
$ws.Range("E20").Value = '  -1.20%  '
$ws.Range("D21").Value = '6.74'
$ws.Range("E21").Value = '  +1.44%  '
$ws.Range("D22").Value = '0.999'
$ws.Range("E22").Value = '  -0.13%  '
$ws.Range("D23").Value = '61.33'
$ws.Range("E23").Value = '  +0.50%  '
$ws.Range("E24").Value = '  +3.90%  '
$ws.Range("D25").Value = '8.62'
$ws.Range("E25").Value = '  +11.18%  '
$ws.Range("D26").Value = '0.994'
$ws.Range("E26").Value = '  +0.28%  '
$ws.Range("D27").Value = '1.33'
$ws.Range("E27").Value = '  +9.10%  '
$ws.Range("E28").Value = '  +0.16%  '
$ws.Range("D29").Value = '167.75'
$ws.Range("E29").Value = '  -2.58%  '
$ws.Range("E30").Value = '  +0.01%  '
$ws.Range("D31").Value = '6.29'
$ws.Range("E31").Value = '  +0.02%  '
$ws.Range("E32").Value = '  +1.07%  '
$ws.Range("E33").Value = '  +0.01%  '
$ws.Range("D34").Value = '1.30'
$ws.Range("E34").Value = '  +2.43%  '
$ws.Range("D35").Value = '0.994'
$ws.Range("E35").Value = '  +1.19%  '
$ws.Range("D36").Value = '0.924'
$ws.Range("E36").Value = '  -3.41%  '
$ws.Range("D37").Value = '4.04'
$ws.Range("E37").Value = '  -0.22%  '
$ws.Range("D38").Value = '1.61'
$ws.Range("E38").Value = '  +5.37%  '
$ws.Range("D39").Value = '38.82'
$ws.Range("E39").Value = '  +3.20%  '
$ws.Range("D40").Value = '150.41'
$ws.Range("E40").Value = '  +7.35%  '
$ws.Range("D41").Value = '0.386'
$ws.Range("E41").Value = '  +0.80%  '
$ws.Range("E42").Value = '  +1.52%  '
$ws.Range("D43").Value = '5.37'
$ws.Range("E43").Value = '  +1.94%  '
$ws.Range("D44").Value = '284.67'
$ws.Range("E44").Value = '  +1.64%  '
$ws.Range("D45").Value = '0.0942'
$ws.Range("E45").Value = '  +1.01%  '
$ws.Range("E46").Value = '  -0.19%  '
$ws.Range("D47").Value = '0.566'
$ws.Range("E47").Value = '  +0.56%  '
$ws.Range("D48").Value = '18.28'
$ws.Range("E48").Value = '  +5.40%  '
$ws.Range("E49").Value = '  +1.38%  '
$ws.Range("B50").Value = 'EnergySwap'
$ws.Range("C50").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D50").Value = '17.75'
$ws.Range("E50").Value = '  +4.63%  '
$ws.Range("B51").Value = 'Polygon'
$ws.Range("C51").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D51").Value = '0.386'
$ws.Range("E51").Value = '  +0.67%  '

# Strip the temporary text-format style so cells end up with the default (no explicit) style,
# matching the original workbook formatting.
$textRange.Style = "Normal"
